$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 33, shifting rows 33:85 down to 34:86.
$ws.Rows(33).Insert()

# Populate the newly inserted row 33 with a new data record (same
# constant columns as every other data row, new measurement values).
$ws.Range("A33").Value = 2
$ws.Range("B33").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C33").Value = "Coquimbo"
$ws.Range("D33").Value = 44671
$ws.Range("E33").Value = 4
$ws.Range("F33").Value = 100112030
$ws.Range("G33").Value = "Poroto granado"
$ws.Range("H33").Value = "Sin especificar"
$ws.Range("I33").Value = "Primera"
$ws.Range("J33").Value = 400
$ws.Range("K33").Value = 17000
$ws.Range("L33").Value = 18000
$ws.Range("M33").Value = 17500
$ws.Range("N33").Value = "$/malla 25 kilos"
$ws.Range("O33").Value = "Provincia de Limarí"
$ws.Range("P33").Value = 700
$ws.Range("Q33").Value = 25
$ws.Range("R33").Value = "Hortaliza"
